$wb = $excel.ActiveWorkbook

function Format-ClaimCell($range) {
    $range.Borders.Item(8).LineStyle = 1
    $range.Borders.Item(9).LineStyle = 1
    $range.Interior.ColorIndex = 2
}

# --- ShipmentInformation sheet ---
$wsShip = $wb.Worksheets.Item("ShipmentInformation")
$wsShip.Range("C2").Value = "PickUp15"
$wsShip.Range("K2").Value = "DropOff445"

# --- Input sheet ---
$wsInput = $wb.Worksheets.Item("Input")

# Row 2
$wsInput.Range("B2").Value = "12-01-2021"
$wsInput.Range("T2").Value = "58327789"
Format-ClaimCell $wsInput.Range("U2")
$wsInput.Range("W2").Value = "FCT915506118178897920"
$wsInput.Range("X2").Value = "FCTEST1004009"
Format-ClaimCell $wsInput.Range("Y2")

# Row 3
$wsInput.Range("B3").Value = "12-01-2021"
$wsInput.Range("T3").Value = "58327790"
Format-ClaimCell $wsInput.Range("U3")
$wsInput.Range("W3").Value = "FCT915518758829686784"
$wsInput.Range("X3").Value = "FCTEST1004010"
Format-ClaimCell $wsInput.Range("Y3")

# --- ClaimDetail sheet ---
$wsClaim = $wb.Worksheets.Item("ClaimDetail")
Format-ClaimCell $wsClaim.Range("C2")
$wsClaim.Range("C3").Value = "Filed"
